# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 4253
$ws1.Range("F3").Value = 2421
$ws1.Range("F10").Value = 130
$ws1.Range("F12").Value = 1585
$ws1.Range("F14").Value = 3269

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 4253
$ws4.Range("F3").Value = 2421
$ws4.Range("F12").Value = 130
$ws4.Range("F16").Value = 1585
$ws4.Range("F18").Value = 3269
